$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Range, $Text) {
    $origStyle = $Range.Style
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "61.069.75"
$ws.Range("E2").Value = "  +1.06%  "
Set-TextValue $ws.Range("D3") "2.360.50"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue $ws.Range("D5") "548.77"
$ws.Range("E5").Value = "  +1.46%  "
Set-TextValue $ws.Range("D6") "139.22"
$ws.Range("E6").Value = "  +2.27%  "
$ws.Range("E7").Value = "  -0.06%  "
Set-TextValue $ws.Range("D8") "0.520"
$ws.Range("E8").Value = "  -0.55%  "
Set-TextValue $ws.Range("D9") "2.361.91"
$ws.Range("E9").Value = "  +1.15%  "
$ws.Range("E10").Value = "  +3.00%  "
$ws.Range("E11").Value = "  +1.80%  "
Set-TextValue $ws.Range("D12") "5.33"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("E13").Value = "  +3.20%  "
Set-TextValue $ws.Range("D14") "25.52"
$ws.Range("E14").Value = "  +4.51%  "
Set-TextValue $ws.Range("D15") "0.0000172"
$ws.Range("E15").Value = "  +7.78%  "
Set-TextValue $ws.Range("D16") "2.786.42"
$ws.Range("E16").Value = "  +1.10%  "
Set-TextValue $ws.Range("D17") "61.312.81"
$ws.Range("E17").Value = "  +1.38%  "
Set-TextValue $ws.Range("D18") "2.358.30"
$ws.Range("E18").Value = "  +0.99%  "
Set-TextValue $ws.Range("D19") "10.98"
$ws.Range("E19").Value = "  +4.63%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Range("D20") "4.14"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue $ws.Range("D21") "320.71"
$ws.Range("E21").Value = "  +1.36%  "
Set-TextValue $ws.Range("D22") "6.61"
$ws.Range("E22").Value = "  +1.02%  "
$ws.Range("E23").Value = "  +0.02%  "
Set-TextValue $ws.Range("D24") "64.11"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E25").Value = "  -7.09%  "
Set-TextValue $ws.Range("D26") "8.85"
$ws.Range("E26").Value = "  +3.20%  "
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D28") "533.32"
$ws.Range("E28").Value = "  +7.01%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D29") "2.473.70"
$ws.Range("E29").Value = "  +0.96%  "
Set-TextValue $ws.Range("D30") "8.21"
$ws.Range("E30").Value = "  +3.75%  "
Set-TextValue $ws.Range("D31") "0.0₃0903"
$ws.Range("E31").Value = "  +1.73%  "
$ws.Range("E32").Value = "  +0.40%  "
Set-TextValue $ws.Range("D33") "0.147"
$ws.Range("E33").Value = "  +1.73%  "
Set-TextValue $ws.Range("D34") "1.84"
$ws.Range("E34").Value = "  +2.87%  "
Set-TextValue $ws.Range("D35") "1.51"
$ws.Range("E35").Value = "  -0.56%  "
Set-TextValue $ws.Range("D36") "0.998"
$ws.Range("E36").Value = "  -0.10%  "
Set-TextValue $ws.Range("D37") "5.61"
$ws.Range("E37").Value = "  +7.75%  "
Set-TextValue $ws.Range("D38") "4.69"
$ws.Range("E38").Value = "  +2.17%  "
Set-TextValue $ws.Range("D39") "1.90"
$ws.Range("E39").Value = "  +6.00%  "
Set-TextValue $ws.Range("D40") "0.380"
$ws.Range("E40").Value = "  +2.20%  "
Set-TextValue $ws.Range("D41") "18.45"
$ws.Range("E41").Value = "  +1.26%  "
Set-TextValue $ws.Range("D42") "144.88"
$ws.Range("E42").Value = "  +5.40%  "
$ws.Range("E43").Value = "  +0.08%  "
Set-TextValue $ws.Range("D44") "41.48"
$ws.Range("E44").Value = "  +3.45%  "
Set-TextValue $ws.Range("D45") "146.92"
$ws.Range("E45").Value = "  +4.37%  "
Set-TextValue $ws.Range("D46") "2.22"
$ws.Range("E46").Value = "  +5.90%  "
Set-TextValue $ws.Range("D47") "3.59"
$ws.Range("E47").Value = "  +1.92%  "
Set-TextValue $ws.Range("D48") "0.0527"
$ws.Range("E48").Value = "  +3.72%  "
Set-TextValue $ws.Range("D49") "19.99"
$ws.Range("E49").Value = "  +3.43%  "
Set-TextValue $ws.Range("D50") "0.578"
$ws.Range("E50").Value = "  +2.21%  "
Set-TextValue $ws.Range("D51") "0.0902"
$ws.Range("E51").Value = "  +0.64%  "
